$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.760.28'
$ws.Range('E2').Value = '  -2.11%  '
$ws.Range('D3').Value = '2.377.86'
$ws.Range('E3').Value = '  -3.85%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.87'
$ws.Range('E5').Value = '  -2.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '87.11'
$ws.Range('E6').Value = '  -5.67%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.532'
$ws.Range('E7').Value = '  -3.57%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.494'
$ws.Range('E9').Value = '  -3.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0845'
$ws.Range('E10').Value = '  -2.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '30.46'
$ws.Range('E11').Value = '  -8.17%  '
$ws.Range('E12').Value = '  -0.51%  '
$ws.Range('D13').Value = '2.738.86'
$ws.Range('E13').Value = '  -4.05%  '
$ws.Range('E14').Value = '  -4.95%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.03'
$ws.Range('E15').Value = '  -3.41%  '
$ws.Range('D16').Value = '2.373.27'
$ws.Range('E16').Value = '  -4.66%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.760'
$ws.Range('E17').Value = '  -3.94%  '
$ws.Range('D18').Value = '40.612.26'
$ws.Range('E18').Value = '  -2.34%  '
$ws.Range('D19').Value = '0.0₃0914'
$ws.Range('E19').Value = '  -3.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.14'
$ws.Range('E20').Value = '  -4.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '68.58'
$ws.Range('E21').Value = '  -3.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.80'
$ws.Range('E22').Value = '  -4.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.63'
$ws.Range('E23').Value = '  -1.94%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.59'
$ws.Range('E24').Value = '  -5.92%  '
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.81'
$ws.Range('E26').Value = '  -7.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.81'
$ws.Range('E27').Value = '  -4.36%  '
$ws.Range('E28').Value = '  -3.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.24'
$ws.Range('E29').Value = '  -4.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.34'
$ws.Range('E30').Value = '  -6.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '154.22'
$ws.Range('E31').Value = '  -1.57%  '
$ws.Range('E32').Value = '  -0.06%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.20'
$ws.Range('E33').Value = '  -4.82%  '
$ws.Range('E34').Value = '  -4.82%  '
$ws.Range('E35').Value = '  -5.74%  '
$ws.Range('E36').Value = '  -2.33%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.79'
$ws.Range('E37').Value = '  -3.53%  '
$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '15.95'
$ws.Range('E38').Value = '  -7.94%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.100'
$ws.Range('E39').Value = '  -3.86%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.71'
$ws.Range('E40').Value = '  -7.53%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.85'
$ws.Range('E41').Value = '  -4.23%  '
$ws.Range('E42').Value = '  -3.83%  '
$ws.Range('D43').Value = '1.965.80'
$ws.Range('E43').Value = '  -1.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0269'
$ws.Range('E44').Value = '  -4.91%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.75'
$ws.Range('E45').Value = '  -5.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.34'
$ws.Range('E46').Value = '  -1.08%  '
$ws.Range('E47').Value = '  -9.21%  '
$ws.Range('D48').Value = '2.598.84'
$ws.Range('E48').Value = '  -4.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '93.25'
$ws.Range('E49').Value = '  -4.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '71.88'
$ws.Range('E50').Value = '  -4.98%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '50.55'
$ws.Range('E51').Value = '  -3.31%  '
